$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Small text fixes (typos)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("La fase avanzad es:", $true, $false, $false, $false, $false, $true, 1, $false, "La fase avanzada es:", 2) | Out-Null

$d.Content.Find.Execute("Iteracion: cuantas", $true, $false, $false, $false, $false, $true, 1, $false, "Iteración: cuantas", 2) | Out-Null

$d.Content.Find.Execute("Estos utilizar la derivada de la función de pérdida para llegar a un mínimo de la pérdida de datos", $true, $false, $false, $false, $false, $true, 1, $false, "Estos utilizan la derivada de la función de pérdida para llegar a un mínimo de la pérdida de datos modificando los pesos de las neuronas.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Re-number the embedded pictures (Shape.Name drives <wp:docPr name>)
#    Shapes appear in document order; map old name -> new name per the
#    commit diff.
# ---------------------------------------------------------------------
$shapeNames = @("image1.png", "image4.png", "image6.png", "image2.png", "image3.png", "image5.png")
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)
    $shp.Name = $shapeNames[$i - 1]
}

# ---------------------------------------------------------------------
# 3) Append the new "neural network notes" paragraphs at the end of the
#    document, right after the existing final sentence, reproducing the
#    same pPr (color 1f1f1f / highlight white) used throughout this
#    section.
# ---------------------------------------------------------------------
$COLOR_1F1F1F = 2039583

function Add-EmptyPara($doc) {
    $r = $doc.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

function Add-TextRun($doc, $text) {
    $r = $doc.Content
    $r.Collapse(0)
    $r.InsertAfter($text)
    $rng = $doc.Range($r.Start, $r.End)
    $rng.Font.Color = $COLOR_1F1F1F
    return @($rng.Start, $rng.End)
}

function Add-BoldTextRun($doc, $text) {
    $r = $doc.Content
    $r.Collapse(0)
    $r.InsertAfter($text)
    $rng = $doc.Range($r.Start, $r.End)
    $rng.Font.Color = $COLOR_1F1F1F
    $rng.Font.Bold = 1
    return @($rng.Start, $rng.End)
}

function Start-NewPara($doc) {
    $r = $doc.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

# -- empty paragraph --
Add-EmptyPara $d

# -- empty paragraph --
Add-EmptyPara $d

# -- "Este maneja como parámetro un **learning_rate**, que es ..." --
Start-NewPara $d
Add-TextRun $d "Este maneja como parámetro un "
Add-BoldTextRun $d "learning_rate"
Add-TextRun $d ", que es el tamaño del paso que tomamos a la hora de buscar el mínimo de la función que estamos utilizando, si usamos un learning rate muy grande podemos saltarnos el mínimo global que estamos buscando en la función que estamos usando, si es muy pequeño puede tardar demasiado en hallar este mínimo."

# -- empty paragraph --
Add-EmptyPara $d

# -- "A la hora de dar el salto ... learning_rate y le multiplicamos ..." --
Start-NewPara $d
Add-TextRun $d "A la hora de dar el salto en la función buscando el minimo, tomamos el learning_rate y le multiplicamos un valor dado por el optimizador que nos dice el paso que daremos en la búsqueda del mínimo de la función."

# -- empty paragraph --
Add-EmptyPara $d

# -- "**Métricas:** son maneras de evaluar ..." --
Start-NewPara $d
Add-BoldTextRun $d "Métricas:"
Add-TextRun $d " son maneras de evaluar la red neuronal en su desempeño mientras se entrena y valida los datos."

# -- empty paragraph --
Add-EmptyPara $d

# -- "**Batch size:** nos dice el tamaño del lote ..." --
Start-NewPara $d
Add-BoldTextRun $d "Batch size:"
Add-TextRun $d " nos dice el tamaño del lote que se va a utilizar para el entrenamiento por cada época"

# -- empty paragraph --
Add-EmptyPara $d

# -- "**Modelo:** maneja función de activación ..." --
Start-NewPara $d
Add-BoldTextRun $d "Modelo:"
Add-TextRun $d " maneja función de activación, función de pérdida y la métrica"

# -- trailing empty paragraph (mirrors the document's original ending) --
Add-EmptyPara $d
